$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P15:AO15 (new optional BioSample attributes)
$headerCols  = @('P', 'Q', 'R', 'S', 'T', 'U', 'V', 'W', 'X', 'Y', 'Z', 'AA', 'AB', 'AC', 'AD', 'AE', 'AF', 'AG', 'AH', 'AI', 'AJ', 'AK', 'AL', 'AM', 'AN', 'AO')
$headerVals  = @('biomaterial_provider', 'birth_date', 'birth_location', 'breeding_history', 'breeding_method', 'cell_line', 'cell_subtype', 'cell_type', 'collected_by', 'collection_date', 'culture_collection', 'death_date', 'disease', 'disease_stage', 'genotype', 'geo_loc_name', 'growth_protocol', 'health_state', 'isolation_source', 'lat_lon', 'phenotype', 'sample_type', 'specimen_voucher', 'store_cond', 'stud_book_number', 'treatment')

for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "15").Value = $headerVals[$i]
}

# Copy the existing optional-field (yellow) header style onto the new header cells
$styleSource = $ws.Range("C15")
$styleSource.Copy()
$ws.Range("P15:AO15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Attach the field-definition comments (popup notes) to the relevant new header cells
$commentRefs  = @('P15', 'U15', 'W15', 'X15', 'Y15', 'Z15', 'AB15', 'AC15', 'AD15', 'AE15', 'AG15', 'AH15', 'AI15', 'AJ15', 'AK15', 'AL15', 'AM15')
$commentTexts = @(
    'name and address of the lab or PI, or a culture collection identifier',
    'Name of the cell line.',
    'Type of cell of the sample or from which the sample was obtained.',
    'Name of persons or institute who collected the sample',
    'Time of sampling (single instance or interval, eg., 2008-01-23T19:23:10, 2008-01-23, 2008-01, 2008, 1952-10-21T11:43Z/1952-10-21T17:43Z, 1952-10-21/1953-02-15, 1952-10/1953-02, 1952/1953)',
    'Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier',
    'list of diseases diagnosed; can include multiple diagnoses. the value of the field depends on host; for humans the terms should be chosen from DO (Disease Ontology), free text for non-human. For DO terms, please see http://gemina.svn.sourceforge.net/viewvc/gemina/trunk/Gemina/ontologies/gemina_symptom.obo?view=log',
    'Stage of disease at the time of sampling.',
    'observed genotype',
    'Geographical origin of the sample; use the appropriate name from the list, http://www.ddbj.nig.ac.jp/sub/country-e.html. Use a colon to separate the country or ocean from more detailed information about the location, eg "Japan:Kanagawa, Hakone, Lake Ashi" ',
    'Health or disease status of sample at time of collection',
    'Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.',
    'The geographical coordinates of the location where the sample was collected. Specify as decimal degrees latitude and longitude in format "d[d.dddd] N|S d[dd.dddd] W|E", eg, 47.94 N 28.12 W',
    'Phenotype of sampled organism. For Phenotypic quality Ontology (PATO) (v1.269) terms, please see http://bioportal.bioontology.org/visualize/44601',
    'Sample type, such as cell culture, mixed culture, tissue sample, whole organism, single cell, metagenomic assembly',
    'Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier',
    'explain how and for how long the soil sample was stored before DNA extraction.'
)

for ($i = 0; $i -lt $commentRefs.Length; $i++) {
    $ws.Range($commentRefs[$i]).AddComment($commentTexts[$i])
}
